$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill the "Max" key-metric column for the first three rows ---
$ws.Range("G2").Value = "Max"
$ws.Range("G3").Value = "Max"
$ws.Range("G4").Value = "Max"

# --- Fill the Models column (new model rows) ---
$ws.Range("B2").Value = "TinyLlama (1.1B)"
$ws.Range("B3").Value = "Phi-2 (2.7B)"
$ws.Range("B4").Value = "Mistral (7B)"
$ws.Range("B5").Value = "LLaMA 3 (8B)"
$ws.Range("B6").Value = "Gemma (2B/7B)"
$ws.Range("B7").Value = "Qwen (1.8B)"

# --- Fill the Quantization column for every model row ---
$ws.Range("D2").Value = "Q4_K_M"
$ws.Range("D3").Value = "Q4_K_M"
$ws.Range("D4").Value = "Q4_K_M"
$ws.Range("D5").Value = "Q4_K_M"
$ws.Range("D6").Value = "Q4_K_M"
$ws.Range("D7").Value = "Q4_K_M"
$ws.Range("D8").Value = "Q4_K_M"
$ws.Range("D9").Value = "Q4_K_M"
$ws.Range("D10").Value = "Q4_K_M"

# --- Re-arrange header row (B1:D1) ---
# Old: B1=Size, C1=Quantization, D1=Models
# New: B1=Models, C1=Size, D1=Quantization
$ws.Range("B1").Value = "Models"
$ws.Range("C1").Value = "Size"
$ws.Range("D1").Value = "Quantization"

# --- Column widths: autofit to new content ---
$ws.Columns.AutoFit()

# --- Selection matches the saved state in the source file ---
$ws.Range("L19").Select()
